# 202506 June full release
# Updates the ACE_landing_page_data workbook:
#  - refreshes the computed metrics on the "ACE_landing_page_data" sheet
#  - renames the ANSP list entry "IAA" -> "UkSATSE" (appended at the end of
#    the alphabetical list, in its own font) and fixes the capitalisation of
#    "Airnav Ireland" -> "AirNav Ireland"
#  - updates the two sheets' saved selections to match where the author left
#    the cursor

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ACE_landing_page_data")
$ws2 = $wb.Worksheets.Item("ANSP")

# ---------------------------------------------------------------------
# 1. Refreshed figures on the data sheet (rows 2-7, columns B-O)
# ---------------------------------------------------------------------
$ws1.Range("B2").Value = 459.51125496142458
$ws1.Range("C2").Value = 9696014411.1017056
$ws1.Range("D2").Value = 21100711.476404805
$ws1.Range("E2").Value = 0.93449148017114469
$ws1.Range("F2").Value = 141.76294716186041
$ws1.Range("G2").Value = 307.81062400875106
$ws1.Range("H2").Value = -0.073924857539470756
$ws1.Range("I2").Value = 0.028658747874637847
$ws1.Range("J2").Value = 0.11077244244084761
$ws1.Range("K2").Value = 0.088513109734926765
$ws1.Range("L2").Value = 0.020863303800442878
$ws1.Range("M2").Value = -0.079620411780101707
$ws1.Range("N2").Value = 98.068454578398587
$ws1.Range("O2").Value = 98.380793802270077
$ws1.Range("B3").Value = 496.1921920725884
$ws1.Range("D3").Value = 18996430.475028183
$ws1.Range("E3").Value = 0.85850273351205741
$ws1.Range("F3").Value = 138.86574885600163
$ws1.Range("G3").Value = 334.43877716159056
$ws1.Range("H3").Value = -0.34633047653699467
$ws1.Range("J3").Value = 0.58283823355771003
$ws1.Range("K3").Value = 0.42204113678545663
$ws1.Range("L3").Value = -0.034209359027104091
$ws1.Range("M3").Value = -0.35798396816584188
$ws1.Range("N3").Value = 95.336237387785417
$ws1.Range("O3").Value = 88.569710629555161
$ws1.Range("B4").Value = 759.08723638187269
$ws1.Range("C4").Value = 9110184227.7325687
$ws1.Range("D4").Value = 12001498.366848476
$ws1.Range("E4").Value = 0.6037116025016791
$ws1.Range("F4").Value = 143.78452530469156
$ws1.Range("G4").Value = 520.91966645465459
$ws1.Range("H4").Value = -0.25298186968904457
$ws1.Range("I4").Value = -0.050059059806503159
$ws1.Range("J4").Value = 0.27164375488192261
$ws1.Range("K4").Value = 0.25218236874053002
$ws1.Range("L4").Value = -0.082298764589214035
$ws1.Range("M4").Value = -0.24633520647259943
$ws1.Range("N4").Value = 92.143188970025463
$ws1.Range("O4").Value = 55.956261828777663
$ws1.Range("B5").Value = 1016.1563763732124
$ws1.Range("D5").Value = 9437783.4364175871
$ws1.Range("E5").Value = 0.48212753794713165
$ws1.Range("F5").Value = 156.67901464721251
$ws1.Range("G5").Value = 691.18216868878562
$ws1.Range("H5").Value = 1.2101255825117465
$ws1.Range("I5").Value = -0.043310217984809829
$ws1.Range("J5").Value = -0.5671332934267298
$ws1.Range("K5").Value = -0.50709841954439916
$ws1.Range("M5").Value = 1.2162320174033208
$ws1.Range("N5").Value = 96.998860741023009
$ws1.Range("O5").Value = 44.003095689306029
$ws1.Range("B6").Value = 459.77313887221686
$ws1.Range("C6").Value = 10024423797.591002
$ws1.Range("D6").Value = 21802978.360545449
$ws1.Range("E6").Value = 0.97814159472057183
$ws1.Range("F6").Value = 144.66761280441656
$ws1.Range("G6").Value = 311.87265740281964
$ws1.Range("H6").Value = -0.0026066897717813564
$ws1.Range("I6").Value = 0.0139008753360228
$ws1.Range("J6").Value = 0.01655070766819855
$ws1.Range("K6").Value = 0.010368208674844448
$ws1.Range("M6").Value = -0.0016746158767168007
$ws1.Range("N6").Value = 101.39008753360228
$ws1.Range("O6").Value = 101.65507076681986
$ws1.Range("B7").Value = 460.97475705648537
$ws1.Range("C7").Value = 9886986037.2383537
$ws1.Range("D7").Value = 21447998.802300699
$ws1.Range("E7").Value = 0.9681040895016485
$ws1.Range("F7").Value = 143.83989558213406
$ws1.Range("G7").Value = 312.3958003699388
$ws1.Range("H7").Value = -0.03564570624302521
$ws1.Range("I7").Value = 0.015666737236474715
$ws1.Range("J7").Value = 0.053209120145662148
$ws1.Range("K7").Value = 0.049247231921604406
$ws1.Range("M7").Value = -0.030183511783170802

# ---------------------------------------------------------------------
# 2. ANSP list maintenance
# ---------------------------------------------------------------------
# Drop the old "IAA" row - everything below shifts up one row.
$ws2.Range("A21").Delete(-4162)

# Re-add it under its new name "UkSATSE" at the bottom of the list, in the
# new default font used for freshly typed entries.
$ws2.Range("A40").Value = "UkSATSE"
$ws2.Range("A40").Font.Name = "Aptos"

# Fix the capitalisation of the Irish ANSP's name.
$ws2.Range("A2").Value = "AirNav Ireland"

# ---------------------------------------------------------------------
# 3. Restore the cursor positions the author left behind
# ---------------------------------------------------------------------
[void]$ws1.Range("B2:O7").Select()

[void]$ws2.Activate()
[void]$ws2.Range("D11").Select()
